$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New "Vendor" column (B) for the existing OutOfDarts line items (rows 8-12)
#    and merge the cell across that block.
# ---------------------------------------------------------------------------
$ws.Range("B8").Value = "OutOfDarts"
$ws.Range("B8:B12").Merge()
$ws.Range("B8:B12").HorizontalAlignment = -4108
$ws.Range("B8:B12").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 2. Row 15 used to be the "245 RPM 12V Gear Motor with Encoder" / AndyMark
#    line. It becomes the Pololu gearmotor line, with a Vendor cell too.
# ---------------------------------------------------------------------------
$ws.Range("E15").Hyperlinks.Delete()
$ws.Range("B15").Value = "Polulu"
$ws.Range("B15").HorizontalAlignment = -4108
$ws.Range("C15").Value = "12V 50:1 Metal Gearmotor with Encoder"
$ws.Range("E15").Value = "https://www.pololu.com/product/4753/resources"
$ws.Range("G15").Value = 51.95
$ws.Range("H15").Formula = "=F15*G15"

# ---------------------------------------------------------------------------
# 3. Row 17 is a brand-new line item: Amazon ball bearing.
# ---------------------------------------------------------------------------
$ws.Range("B17").Value = "Amazon"
$ws.Range("B17").HorizontalAlignment = -4108
$ws.Range("C17").Value = "Deep Groove Ball Bearing 80x100x10mm"
$ws.Range("E17").Value = "https://www.amazon.com/uxcell-6816-2RS-Bearing-80x100x10mm-Bearings/dp/B07RQ4RXDR/ref=sr_1_2_sspa?crid=245V7BXKB1Z1W&keywords=large%2Bbearing&qid=1706830027&sprefix=large%2Bbearing%2Caps%2C144&sr=8-2-spons&sp_csd=d2lkZ2V0TmFtZT1zcF9hdGY&th=1"
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 13.99
$ws.Range("H17").Formula = "=F17*G17"
$ws.Range("G17").Style = "Currency"
$ws.Range("H17").Style = "Currency"
$ws.Range("G17").NumberFormat = """$""#,##0.00"
$ws.Range("H17").NumberFormat = """$""#,##0.00"

# ---------------------------------------------------------------------------
# Re-create every hyperlink (the ones that survived unchanged, plus the new
# Pololu / Amazon ones) so the worksheet's hyperlink table matches the final
# layout: E8, E9, E10, E11, E12, E17.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E8"), "https://outofdarts.com/products/loki-130-3s-high-rpm-neo-motor-for-nerf-blasters") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E9"), "https://outofdarts.com/products/worker-10-round-talon-short-dart-magazine?variant=39475795984426") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E10"), "https://outofdarts.com/products/nightingale-flywheel-pair") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E11"), "https://outofdarts.com/products/n20-metal-gear-motor-micro-size-300-3000rmp-multiple-options") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E12"), "https://outofdarts.com/products/micro-switch-1a-for-mosfet-builds") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E17"), "https://www.pololu.com/product/4753/resources") | Out-Null

$ws.Range("E8:E12").Style = "Hyperlink"
$ws.Range("E15").Style = "Hyperlink"
$ws.Range("E17").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 4. Total formula range / label stay put (SUM still covers H8:H23) - Excel
#    recalculates the value automatically once the new rows have numbers.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 5. Selection moved to J2 in the saved file.
# ---------------------------------------------------------------------------
$ws.Range("J2").Select() | Out-Null

$wb.Save()
